$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update goal values in column B (experiment A length / experiment B .xosc correction)
$ws.Range("B2").Value = 88
$ws.Range("B3").Value = -259
$ws.Range("B5").Value = 0.343
$ws.Range("B6").Value = -0.618
$ws.Range("B7").Value = 0.343
$ws.Range("B8").Value = -0.618

# Move / update the active selection to B9 (matches the saved cursor position)
$ws.Range("B9").Select()
